# edit.ps1 - apply resume content updates to $word.ActiveDocument
#
# Summary of changes (see commit message / diff):
#  1. "21 years of experience" -> "15+ years of experience" (PROFESSIONAL SUMMARY)
#  2. FLEEM bullet (RESEARCH DIRECTOR section) - expand Twilio API description
#  3. Salsa Labs CRM bullet - append "used by tens of thousands of users simultaneously"
#  4. Salsa Labs mapping bullet - append "interfacing with Government and Activism APIs"
#  5. NEW bullet after Salsa Labs "Collaborated with political strategists..." line
#  6. NEW bullet after Praxis Project "Managed technology infrastructure..." line
#  7. NEW bullet after Lake Research Partners "Developed innovative approaches..." line
#  8. NEW bullet after Feldman Group "Enhanced value of research deliverables..." line

$d = $word.ActiveDocument

function Replace-InParagraph($index, $findText, $replaceText) {
    $para = $d.Paragraphs.Item($index)
    $para.Range.Find.Execute($findText, $true, $false, $false, $false, $false, `
        $true, 1, $false, $replaceText, 2) | Out-Null
}

function Insert-BulletAfter($index, $bulletText) {
    $para = $d.Paragraphs.Item($index)
    $para.Range.InsertParagraphAfter() | Out-Null
    $d.Paragraphs.Item($index + 1).Range.Text = $bulletText
}

# Process insertions from the bottom of the document upward so that
# paragraph indices for not-yet-processed items remain valid.

# 8. Feldman Group - add PHP/MySQL training bullet
Insert-BulletAfter 73 "• Trained staff on PHP/MySQL for data analysis and reporting systems"

# 7. Lake Research Partners - add Python tooling training bullet
Insert-BulletAfter 65 "• Trained staff on building Python tooling for report generation and analysis"

# 6. Praxis Project - add Drupal sites bullet
Insert-BulletAfter 57 "• Architected and developed 25 Drupal sites to integrate with membership databases, activism CRMs and government agencies, under guidelines from Kellogg Foundation and Robert Wood Johnson Foundation"

# 5. Salsa Labs - add billions of records bullet
Insert-BulletAfter 49 "• Handled billions of records with millions of columns in high-performance CRM system"

# 4. Salsa Labs - mapping/visualization bullet gets API integration detail
Replace-InParagraph 48 "Integrated mapping and visualization tools for political campaign data analysis" `
    "Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs"

# 3. Salsa Labs - Java-based CRM bullet gets scale detail
Replace-InParagraph 45 "Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system" `
    "Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands of users simultaneously"

# 2. RESEARCH DIRECTOR - FLEEM bullet gets expanded Twilio description
Replace-InParagraph 37 "Twilio API for thousands of simultaneous phone calls" `
    "Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys"

# 1. PROFESSIONAL SUMMARY - years of experience
Replace-InParagraph 4 "21 years of experience" "15+ years of experience"
